$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93-108 down to 94-109.
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with the new data record.
$ws.Range("A93").Value = 3
$ws.Range("B93").Value = "Femacal de La Calera"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44474
$ws.Range("E93").Value = 5
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100101
$ws.Range("H93").Value = "Berries"
$ws.Range("I93").Value = 100101001
$ws.Range("J93").Value = "Arándano (blue)"
$ws.Range("K93").Value = "Sin especificar"
$ws.Range("L93").Value = "Primera"
$ws.Range("M93").Value = 45
$ws.Range("N93").Value = 10000
$ws.Range("O93").Value = 10000
$ws.Range("P93").Value = 10000
$ws.Range("Q93").Value = '$/bandeja 2 kilos'
$ws.Range("R93").Value = "Provincia de Quillota"
$ws.Range("S93").Value = 5000
$ws.Range("T93").Value = 2
